$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new category "Laboratório" as a new row at the bottom of column A,
# copying the style of the existing data rows (e.g. A13).
$ws.Range("A13").Copy()
$ws.Range("A14").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A14").Value = "Laboratório"

# Move selection to the next empty cell below, as Excel would after data entry.
$ws.Range("A15").Select()
